$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Keep the first sheet active/selected (matches original workbook view state)
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
